$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived data: full 3x3 cross product of clusters (ECs, FAPs, MuSCs)
# as sending cluster x target cluster, replacing the old data that included
# "Neutrophils" and only had partial combinations.

$data = @(
  @("ECs",  "Gdnf", "Gfra1", "ECs",   2, 1,                   0.4995885, 0.999177,  0.08073427759634712, 0.07892927665008602, 2, 1, 0.452919,          0.905838,          0.03428003430836644, 0.03049201402277582, 0.2262731238315,   0.905092495326,    0.002767573805863959, 0.002406712610421975),
  @("ECs",  "Gdnf", "Gfra1", "FAPs",  2, 1,                   0.4995885, 0.999177,  0.08073427759634712, 0.07892927665008602, 3, 1, 3.282732333333333, 9.848196999999999, 0.2484598283839946,  0.3315066943791922,  1.6400153223115,   9.840091933868999, 0.02005922475629419,  0.02616558359201078),
  @("ECs",  "Gdnf", "Gfra1", "MuSCs", 2, 1,                   0.4995885, 0.999177,  0.08073427759634712, 0.07892927665008602, 2, 1, 9.476675,          18.95335,          0.7172601373076389,  0.638001291598032,   4.7344378482375,   18.93775139295,    0.05790747903418897, 0.05035698044765326),
  @("FAPs", "Gdnf", "Gfra1", "ECs",   2, 0.6666666666666666,  0.2830243333333333, 0.849073, 0.0457371718801662,  0.06707191790155147, 2, 1, 0.452919,          0.905838,          0.03428003430836644, 0.03049201402277582, 0.128187098029,    0.769122588174,    0.00156787182121975,  0.002045157861188576),
  @("FAPs", "Gdnf", "Gfra1", "FAPs",  2, 0.6666666666666666,  0.2830243333333333, 0.849073, 0.0457371718801662,  0.06707191790155147, 3, 1, 3.282732333333333, 9.848196999999999, 0.2484598283839946,  0.3315066943791922,  0.9290931301534443, 8.361838171380999, 0.01136384987611536,  0.0222347897892159),
  @("FAPs", "Gdnf", "Gfra1", "MuSCs", 2, 0.6666666666666666,  0.2830243333333333, 0.849073, 0.0457371718801662,  0.06707191790155147, 2, 1, 9.476675,          18.95335,          0.7172601373076389,  0.638001291598032,   2.682129624091667, 16.09277774455,    0.03280545018283109,  0.042791970251147),
  @("MuSCs","Gdnf", "Gfra1", "ECs",   2, 1,                   5.4054465, 10.810893, 0.8735285505234867,  0.8539988054483625,  2, 1, 0.452919,          0.905838,          0.03428003430836644, 0.03049201402277582, 2.4482294233335,   9.792917693334001, 0.02994458868128273,  0.02604014355116527),
  @("MuSCs","Gdnf", "Gfra1", "FAPs",  2, 1,                   5.4054465, 10.810893, 0.8735285505234867,  0.8539988054483625,  3, 1, 3.282732333333333, 9.848196999999999, 0.2484598283839946,  0.3315066943791922,  17.7446340016535,  106.467804009921,  0.2170367537515851,   0.2831063209979656),
  @("MuSCs","Gdnf", "Gfra1", "MuSCs", 2, 1,                   5.4054465, 10.810893, 0.8735285505234867,  0.8539988054483625,  2, 1, 9.476675,          18.95335,          0.7172601373076389,  0.638001291598032,   51.2256597103875,  204.90263884155,   0.6265472080906188,   0.5448523408992317)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $vals = $data[$i]
  for ($j = 0; $j -lt $vals.Count; $j++) {
    $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
  }
}
